$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- New "Index" column R (with Q as a red/status spacer like O) ---
$ws.Range("R5").Value = "Create Unit Test for class"
$ws.Range("R6").Value = "Loading icon"
$ws.Range("R7").Value = "Make ppt"
$ws.Range("R8").Value = "Design"
$ws.Range("R9").Value = "Picture rep"
$ws.Range("R10").Value = "Logo?"
$ws.Range("R11").Value = "Big merge"

# Style Q5:Q11 like the existing status spacer column O (style index 5 - amber fill)
$ws.Range("Q5:Q11").Interior.Color = $ws.Range("O6").Interior.Color

# R7 and R8 use the bold-font style like P6/P8 (style index 9)
$ws.Range("R7").Font.Bold = $true
$ws.Range("R8").Font.Bold = $true

# --- Status updates in column O (amber -> green) ---
$ws.Range("O5").Interior.Color = $ws.Range("O7").Interior.Color
$ws.Range("O9").Interior.Color = $ws.Range("O7").Interior.Color

# --- Remove now-unused bold style from T4 / P10 ---
$ws.Range("T4").Font.Bold = $false
$ws.Range("P10").Font.Bold = $false

# --- Column R width to match column P ---
$ws.Columns.Item(18).ColumnWidth = $ws.Columns.Item(16).ColumnWidth

# --- Window/selection state ---
$ws.Range("R11").Select()

$excel.ActiveWindow.WindowState = -4143
$excel.Width = 19380
$excel.Height = 10260
